$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 2 (matches the author's recorded selection before deleting),
# then delete the entire row. This removes the "setup00 / Software
# installation" row, shifting all subsequent rows up by one and
# adjusting the relative formula references automatically.
$ws.Range("A2:XFD2").Select()
$ws.Rows("2:2").Delete()
